$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.091.50"
$ws.Range("E2").Value = "  +1.61%  "
$ws.Range("D3").Value = "2.248.79"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'318.38"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").Value = "'100.63"
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("D7").Value = "'0.574"
$ws.Range("E7").Value = "  -1.43%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'0.544"
$ws.Range("E9").Value = "  -3.30%  "
$ws.Range("D10").Value = "'36.75"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").Value = "'0.0826"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "'7.51"
$ws.Range("E12").Value = "  -2.77%  "
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D14").Value = "2.591.60"
$ws.Range("E14").Value = "  +0.58%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.271.48"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.850"
$ws.Range("E16").Value = "  -1.91%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'14.20"
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("D18").Value = "43.983.20"
$ws.Range("E18").Value = "  +1.56%  "
$ws.Range("D19").Value = "'13.54"
$ws.Range("E19").Value = "  -4.29%  "
$ws.Range("D20").Value = "0.0₃0974"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "'6.45"
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("D22").Value = "'65.35"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'3.09"
$ws.Range("E23").Value = "  -4.11%  "
$ws.Range("D24").Value = "'234.39"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").Value = "'2.06"
$ws.Range("E25").Value = "  -5.82%  "
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").Value = "'10.49"
$ws.Range("E27").Value = "  +4.03%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "'38.25"
$ws.Range("E28").Value = "  +3.90%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.21"
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("D30").Value = "'6.08"
$ws.Range("E30").Value = "  -5.58%  "
$ws.Range("D31").Value = "'158.38"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("D32").Value = "'20.09"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("D33").Value = "'0.0847"
$ws.Range("E33").Value = "  -3.86%  "
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("D35").Value = "'3.23"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'0.113"
$ws.Range("E36").Value = "  +7.76%  "
$ws.Range("D37").Value = "'1.95"
$ws.Range("E37").Value = "  +4.35%  "
$ws.Range("D38").Value = "'0.118"
$ws.Range("E38").Value = "  -2.35%  "
$ws.Range("D39").Value = "'16.28"
$ws.Range("E39").Value = "  +13.84%  "
$ws.Range("D40").Value = "'3.66"
$ws.Range("E40").Value = "  -1.84%  "
$ws.Range("D41").Value = "'4.15"
$ws.Range("E41").Value = "  -6.33%  "
$ws.Range("D42").Value = "'0.0313"
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "1.768.35"
$ws.Range("E44").Value = "  -2.22%  "
$ws.Range("D45").Value = "'74.94"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").Value = "'0.195"
$ws.Range("E46").Value = "  -3.92%  "
$ws.Range("D47").Value = "'81.22"
$ws.Range("E47").Value = "  -4.17%  "
$ws.Range("D48").Value = "'5.15"
$ws.Range("E48").Value = "  -2.78%  "
$ws.Range("D49").Value = "'102.91"
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.65"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'57.43"
$ws.Range("E51").Value = "  -1.95%  "